$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column BK (63) so the existing "Township" column
# (and its LEGDAT tag above it) shifts right to BL, and the new column takes
# over BK for "Model Group" (tagged as COMDAT, like the other commercial
# building fields).
$ws.Columns.Item(63).Insert()

$ws.Cells.Item(1, 63).Value = "COMDAT"
$ws.Cells.Item(2, 63).Value = "Model Group"
